# Auto-generated edit script: update Leve profit-tracking values
# (currentAveragePrice / Nq / Hq, LevePrice Nq/Hq, LeveProfit Nq/Hq)
# per sheet, reflecting a refreshed market-board data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3640.0588  # H64: 3782.2942 -> 3640.0588
$ws.Cells.Item(64, 9).Value = 3241.1667  # I64: 3667.75 -> 3241.1667
$ws.Cells.Item(64, 10).Value = 4597.4  # J64: 3884.111 -> 4597.4
$ws.Cells.Item(64, 11).Value = 3241.1667  # K64: 3667.75 -> 3241.1667
$ws.Cells.Item(64, 12).Value = 4597.4  # L64: 3884.111 -> 4597.4
$ws.Cells.Item(64, 13).Value = -2993.1667  # M64: -3419.75 -> -2993.1667
$ws.Cells.Item(64, 14).Value = -5093.4  # N64: -4380.111 -> -5093.4
$ws.Cells.Item(67, 8).Value = 3640.0588  # H67: 3782.2942 -> 3640.0588
$ws.Cells.Item(67, 9).Value = 3241.1667  # I67: 3667.75 -> 3241.1667
$ws.Cells.Item(67, 10).Value = 4597.4  # J67: 3884.111 -> 4597.4
$ws.Cells.Item(67, 11).Value = 3241.1667  # K67: 3667.75 -> 3241.1667
$ws.Cells.Item(67, 12).Value = 4597.4  # L67: 3884.111 -> 4597.4
$ws.Cells.Item(67, 13).Value = -2383.1667  # M67: -2809.75 -> -2383.1667
$ws.Cells.Item(67, 14).Value = -6313.4  # N67: -5600.111 -> -6313.4
$ws.Cells.Item(107, 8).Value = 787.0526  # H107: 1488.875 -> 787.0526
$ws.Cells.Item(107, 9).Value = 775.2222  # I107: 1448.625 -> 775.2222
$ws.Cells.Item(107, 10).Value = 1000  # J107: 1569.375 -> 1000
$ws.Cells.Item(107, 11).Value = 775.2222  # K107: 1448.625 -> 775.2222
$ws.Cells.Item(107, 12).Value = 1000  # L107: 1569.375 -> 1000
$ws.Cells.Item(107, 13).Value = 1144.7778  # M107: 471.375 -> 1144.7778
$ws.Cells.Item(107, 14).Value = -4840  # N107: -5409.375 -> -4840
$ws.Cells.Item(131, 8).Value = 3603.3333  # H131: 2652.2222 -> 3603.3333
$ws.Cells.Item(131, 9).Value = 1933.75  # I131: 704.6875 -> 1933.75
$ws.Cells.Item(131, 10).Value = 4630.769  # J131: 7446.154 -> 4630.769
$ws.Cells.Item(131, 11).Value = 5801.25  # K131: 2114.0625 -> 5801.25
$ws.Cells.Item(131, 12).Value = 13892.307  # L131: 22338.462 -> 13892.307
$ws.Cells.Item(131, 13).Value = -761.25  # M131: 2925.9375 -> -761.25
$ws.Cells.Item(131, 14).Value = -23972.307  # N131: -32418.462 -> -23972.307
$ws.Cells.Item(135, 8).Value = 3620.8438  # H135: 3566.8386 -> 3620.8438
$ws.Cells.Item(135, 9).Value = 2979.6086  # I135: 2590.5715 -> 2979.6086
$ws.Cells.Item(135, 10).Value = 5259.5557  # J135: 12678.667 -> 5259.5557
$ws.Cells.Item(135, 11).Value = 26816.4774  # K135: 23315.1435 -> 26816.4774
$ws.Cells.Item(135, 12).Value = 47336.0013  # L135: 114108.003 -> 47336.0013
$ws.Cells.Item(135, 13).Value = -24281.4774  # M135: -20780.1435 -> -24281.4774
$ws.Cells.Item(135, 14).Value = -52406.0013  # N135: -119178.003 -> -52406.0013
$ws.Cells.Item(137, 8).Value = 1890.8334  # H137: 1566 -> 1890.8334
$ws.Cells.Item(137, 9).Value = 1442.1666  # I137: 1172.4706 -> 1442.1666
$ws.Cells.Item(137, 10).Value = 2115.1667  # J137: 2123.5 -> 2115.1667
$ws.Cells.Item(137, 11).Value = 4326.4998  # K137: 3517.4118 -> 4326.4998
$ws.Cells.Item(137, 12).Value = 6345.500100000001  # L137: 6370.5 -> 6345.500100000001
$ws.Cells.Item(137, 13).Value = -1776.4998  # M137: -967.4118000000003 -> -1776.4998
$ws.Cells.Item(137, 14).Value = -11445.5001  # N137: -11470.5 -> -11445.5001
$ws.Cells.Item(138, 8).Value = 3842.6309  # H138: 3986.0483 -> 3842.6309
$ws.Cells.Item(138, 9).Value = 2674.7144  # I138: 2840.4736 -> 2674.7144
$ws.Cells.Item(138, 10).Value = 4400.0454  # J138: 4492.2324 -> 4400.0454
$ws.Cells.Item(138, 11).Value = 8024.1432  # K138: 8521.4208 -> 8024.1432
$ws.Cells.Item(138, 12).Value = 13200.1362  # L138: 13476.6972 -> 13200.1362
$ws.Cells.Item(138, 13).Value = -2884.1432  # M138: -3381.4208 -> -2884.1432
$ws.Cells.Item(138, 14).Value = -23480.1362  # N138: -23756.6972 -> -23480.1362
$ws.Cells.Item(139, 8).Value = 78020  # H139: 78050 -> 78020
$ws.Cells.Item(139, 10).Value = 78020  # J139: 78050 -> 78020
$ws.Cells.Item(139, 12).Value = 78020  # L139: 78050 -> 78020
$ws.Cells.Item(139, 14).Value = -88300  # N139: -88330 -> -88300
$ws.Cells.Item(140, 8).Value = 97950  # H140: 98000 -> 97950
$ws.Cells.Item(140, 10).Value = 97950  # J140: 98000 -> 97950
$ws.Cells.Item(140, 12).Value = 97950  # L140: 98000 -> 97950
$ws.Cells.Item(140, 14).Value = -108310  # N140: -108360 -> -108310
$ws.Cells.Item(141, 8).Value = 4085.476  # H141: 4438.6113 -> 4085.476
$ws.Cells.Item(141, 9).Value = 4343.125  # I141: 4599.3335 -> 4343.125
$ws.Cells.Item(141, 10).Value = 3261  # J141: 3635 -> 3261
$ws.Cells.Item(141, 11).Value = 13029.375  # K141: 13798.0005 -> 13029.375
$ws.Cells.Item(141, 12).Value = 9783  # L141: 10905 -> 9783
$ws.Cells.Item(141, 13).Value = -7849.375  # M141: -8618.000499999998 -> -7849.375
$ws.Cells.Item(141, 14).Value = -20143  # N141: -21265 -> -20143

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 816486.6  # H32: 956919.0600000001 -> 816486.6
$ws.Cells.Item(32, 9).Value = 10458.6875  # I32: 11190.881 -> 10458.6875
$ws.Cells.Item(32, 10).Value = 7264710  # J32: 10887065 -> 7264710
$ws.Cells.Item(32, 11).Value = 10458.6875  # K32: 11190.881 -> 10458.6875
$ws.Cells.Item(32, 12).Value = 7264710  # L32: 10887065 -> 7264710
$ws.Cells.Item(32, 13).Value = -10171.6875  # M32: -10903.881 -> -10171.6875
$ws.Cells.Item(32, 14).Value = -7265284  # N32: -10887639 -> -7265284
$ws.Cells.Item(61, 8).Value = 4800.4  # H61: 25993.1 -> 4800.4
$ws.Cells.Item(61, 10).Value = 3057  # J61: 109020.5 -> 3057
$ws.Cells.Item(61, 12).Value = 3057  # L61: 109020.5 -> 3057
$ws.Cells.Item(61, 14).Value = -3481  # N61: -109444.5 -> -3481
$ws.Cells.Item(122, 8).Value = 3494.762  # H122: 16995.273 -> 3494.762
$ws.Cells.Item(122, 9).Value = 3419.0789  # I122: 18273.426 -> 3419.0789
$ws.Cells.Item(122, 11).Value = 10257.2367  # K122: 54820.278 -> 10257.2367
$ws.Cells.Item(122, 13).Value = -7807.236699999999  # M122: -52370.278 -> -7807.236699999999
$ws.Cells.Item(136, 8).Value = 4800.4  # H136: 25993.1 -> 4800.4
$ws.Cells.Item(136, 10).Value = 3057  # J136: 109020.5 -> 3057
$ws.Cells.Item(136, 12).Value = 9171  # L136: 327061.5 -> 9171
$ws.Cells.Item(136, 14).Value = -14271  # N136: -332161.5 -> -14271

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(38, 8).Value = 9000  # H38: 0 -> 9000
$ws.Cells.Item(38, 10).Value = 9000  # J38: 0 -> 9000
$ws.Cells.Item(38, 12).Value = 9000  # L38: 0 -> 9000
$ws.Cells.Item(38, 14).Value = -9832  # N38: None -> -9832
$ws.Cells.Item(134, 8).Value = 5616.759  # H134: 6001.4814 -> 5616.759
$ws.Cells.Item(134, 9).Value = 570.2727  # I134: 604.7143 -> 570.2727
$ws.Cells.Item(134, 10).Value = 21477.143  # J134: 24890.166 -> 21477.143
$ws.Cells.Item(134, 11).Value = 1710.8181  # K134: 1814.1429 -> 1710.8181
$ws.Cells.Item(134, 12).Value = 64431.429  # L134: 74670.49800000001 -> 64431.429
$ws.Cells.Item(134, 13).Value = 824.1819  # M134: 720.8571000000002 -> 824.1819
$ws.Cells.Item(134, 14).Value = -69501.429  # N134: -79740.49800000001 -> -69501.429

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(25, 8).Value = 2432.6  # H25: 111 -> 2432.6
$ws.Cells.Item(25, 10).Value = 3013  # J25: 0 -> 3013
$ws.Cells.Item(25, 12).Value = 3013  # L25: 0 -> 3013
$ws.Cells.Item(25, 14).Value = -3361  # N25: None -> -3361
$ws.Cells.Item(62, 8).Value = 4030.875  # H62: 4053.9565 -> 4030.875
$ws.Cells.Item(62, 9).Value = 3930.4546  # I62: 3975.9092 -> 3930.4546
$ws.Cells.Item(62, 10).Value = 4115.846  # J62: 4125.5 -> 4115.846
$ws.Cells.Item(62, 11).Value = 3930.4546  # K62: 3975.9092 -> 3930.4546
$ws.Cells.Item(62, 12).Value = 4115.846  # L62: 4125.5 -> 4115.846
$ws.Cells.Item(62, 13).Value = -3306.4546  # M62: -3351.9092 -> -3306.4546
$ws.Cells.Item(62, 14).Value = -5363.846  # N62: -5373.5 -> -5363.846
$ws.Cells.Item(65, 8).Value = 4030.875  # H65: 4053.9565 -> 4030.875
$ws.Cells.Item(65, 9).Value = 3930.4546  # I65: 3975.9092 -> 3930.4546
$ws.Cells.Item(65, 10).Value = 4115.846  # J65: 4125.5 -> 4115.846
$ws.Cells.Item(65, 11).Value = 19652.273  # K65: 19879.546 -> 19652.273
$ws.Cells.Item(65, 12).Value = 20579.23  # L65: 20627.5 -> 20579.23
$ws.Cells.Item(65, 13).Value = -16532.273  # M65: -16759.546 -> -16532.273
$ws.Cells.Item(65, 14).Value = -26819.23  # N65: -26867.5 -> -26819.23
$ws.Cells.Item(132, 8).Value = 2256.353  # H132: 2762.5386 -> 2256.353
$ws.Cells.Item(132, 9).Value = 1049.6364  # I132: 1500.1666 -> 1049.6364
$ws.Cells.Item(132, 10).Value = 4468.6665  # J132: 3844.5715 -> 4468.6665
$ws.Cells.Item(132, 11).Value = 3148.9092  # K132: 4500.4998 -> 3148.9092
$ws.Cells.Item(132, 12).Value = 13405.9995  # L132: 11533.7145 -> 13405.9995
$ws.Cells.Item(132, 13).Value = -618.9092000000001  # M132: -1970.4998 -> -618.9092000000001
$ws.Cells.Item(132, 14).Value = -18465.9995  # N132: -16593.7145 -> -18465.9995
$ws.Cells.Item(138, 8).Value = 49960  # H138: 50000 -> 49960
$ws.Cells.Item(138, 10).Value = 49960  # J138: 50000 -> 49960
$ws.Cells.Item(138, 12).Value = 49960  # L138: 50000 -> 49960
$ws.Cells.Item(138, 14).Value = -60240  # N138: -60280 -> -60240

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 16667899  # H131: 12501258 -> 16667899
$ws.Cells.Item(131, 10).Value = 18519464  # J131: 13514564 -> 18519464
$ws.Cells.Item(131, 12).Value = 55558392  # L131: 40543692 -> 55558392
$ws.Cells.Item(131, 14).Value = -55568472  # N131: -40553772 -> -55568472

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 5791213  # H113: 6472347.5 -> 5791213
$ws.Cells.Item(113, 9).Value = 16668228  # I113: 20001474 -> 16668228
$ws.Cells.Item(113, 10).Value = 771051.9399999999  # J113: 835211.8 -> 771051.9399999999
$ws.Cells.Item(113, 11).Value = 16668228  # K113: 20001474 -> 16668228
$ws.Cells.Item(113, 12).Value = 771051.9399999999  # L113: 835211.8 -> 771051.9399999999
$ws.Cells.Item(113, 13).Value = -16666058  # M113: -19999304 -> -16666058
$ws.Cells.Item(113, 14).Value = -775391.9399999999  # N113: -839551.8 -> -775391.9399999999
$ws.Cells.Item(122, 8).Value = 2553.7273  # H122: 3258.2856 -> 2553.7273
$ws.Cells.Item(122, 9).Value = 1946.25  # I122: 2500 -> 1946.25
$ws.Cells.Item(122, 10).Value = 2900.8572  # J122: 3561.6 -> 2900.8572
$ws.Cells.Item(122, 11).Value = 5838.75  # K122: 7500 -> 5838.75
$ws.Cells.Item(122, 12).Value = 8702.571599999999  # L122: 10684.8 -> 8702.571599999999
$ws.Cells.Item(122, 13).Value = -3388.75  # M122: -5050 -> -3388.75
$ws.Cells.Item(122, 14).Value = -13602.5716  # N122: -15584.8 -> -13602.5716
$ws.Cells.Item(126, 8).Value = 9805390  # H126: 16668125 -> 9805390
$ws.Cells.Item(126, 9).Value = 1255.625  # I126: 1345.6666 -> 1255.625
$ws.Cells.Item(126, 10).Value = 18520176  # J126: 41668292 -> 18520176
$ws.Cells.Item(126, 11).Value = 3766.875  # K126: 4036.9998 -> 3766.875
$ws.Cells.Item(126, 12).Value = 55560528  # L126: 125004876 -> 55560528
$ws.Cells.Item(126, 13).Value = -1296.875  # M126: -1566.9998 -> -1296.875
$ws.Cells.Item(126, 14).Value = -55565468  # N126: -125009816 -> -55565468
$ws.Cells.Item(138, 8).Value = 69850  # H138: 0 -> 69850
$ws.Cells.Item(138, 10).Value = 69850  # J138: 0 -> 69850
$ws.Cells.Item(138, 12).Value = 69850  # L138: 0 -> 69850
$ws.Cells.Item(138, 14).Value = -80130  # N138: None -> -80130
$ws.Cells.Item(139, 8).Value = 54797.25  # H139: 52200.363 -> 54797.25
$ws.Cells.Item(139, 10).Value = 54797.25  # J139: 52200.363 -> 54797.25
$ws.Cells.Item(139, 12).Value = 54797.25  # L139: 52200.363 -> 54797.25
$ws.Cells.Item(139, 14).Value = -65077.25  # N139: -62480.363 -> -65077.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4025.762  # H7: 3881.7917 -> 4025.762
$ws.Cells.Item(7, 9).Value = 3226.6667  # I7: 3249.2222 -> 3226.6667
$ws.Cells.Item(7, 10).Value = 4345.4  # J7: 4261.3335 -> 4345.4
$ws.Cells.Item(7, 11).Value = 3226.6667  # K7: 3249.2222 -> 3226.6667
$ws.Cells.Item(7, 12).Value = 4345.4  # L7: 4261.3335 -> 4345.4
$ws.Cells.Item(7, 13).Value = -3114.6667  # M7: -3137.2222 -> -3114.6667
$ws.Cells.Item(7, 14).Value = -4569.4  # N7: -4485.3335 -> -4569.4
$ws.Cells.Item(32, 8).Value = 1000  # H32: 927.5 -> 1000
$ws.Cells.Item(32, 9).Value = 1000  # I32: 927.5 -> 1000
$ws.Cells.Item(32, 11).Value = 1000  # K32: 927.5 -> 1000
$ws.Cells.Item(32, 13).Value = -683  # M32: -610.5 -> -683
$ws.Cells.Item(40, 8).Value = 22225304  # H40: 22225234 -> 22225304
$ws.Cells.Item(40, 9).Value = 3133.3333  # I40: 2996 -> 3133.3333
$ws.Cells.Item(40, 10).Value = 33336388  # J40: 40003024 -> 33336388
$ws.Cells.Item(40, 11).Value = 3133.3333  # K40: 2996 -> 3133.3333
$ws.Cells.Item(40, 12).Value = 33336388  # L40: 40003024 -> 33336388
$ws.Cells.Item(40, 13).Value = -2997.3333  # M40: -2860 -> -2997.3333
$ws.Cells.Item(40, 14).Value = -33336660  # N40: -40003296 -> -33336660
$ws.Cells.Item(126, 8).Value = 4025.762  # H126: 3881.7917 -> 4025.762
$ws.Cells.Item(126, 9).Value = 3226.6667  # I126: 3249.2222 -> 3226.6667
$ws.Cells.Item(126, 10).Value = 4345.4  # J126: 4261.3335 -> 4345.4
$ws.Cells.Item(126, 11).Value = 9680.000100000001  # K126: 9747.6666 -> 9680.000100000001
$ws.Cells.Item(126, 12).Value = 13036.2  # L126: 12784.0005 -> 13036.2
$ws.Cells.Item(126, 13).Value = -7210.000100000001  # M126: -7277.6666 -> -7210.000100000001
$ws.Cells.Item(126, 14).Value = -17976.2  # N126: -17724.0005 -> -17976.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 0  # H62: 5732.6665 -> 0
$ws.Cells.Item(62, 9).Value = 0  # I62: 5099 -> 0
$ws.Cells.Item(62, 10).Value = 0  # J62: 7000 -> 0
$ws.Cells.Item(62, 11).Value = 0  # K62: 5099 -> 0
$ws.Cells.Item(62, 12).Value = 0  # L62: 7000 -> 0
$ws.Cells.Item(62, 13).ClearContents()  # M62: -4475 -> (cleared)
$ws.Cells.Item(62, 14).ClearContents()  # N62: -8248 -> (cleared)
$ws.Cells.Item(65, 8).Value = 0  # H65: 5732.6665 -> 0
$ws.Cells.Item(65, 9).Value = 0  # I65: 5099 -> 0
$ws.Cells.Item(65, 10).Value = 0  # J65: 7000 -> 0
$ws.Cells.Item(65, 11).Value = 0  # K65: 25495 -> 0
$ws.Cells.Item(65, 12).Value = 0  # L65: 35000 -> 0
$ws.Cells.Item(65, 13).ClearContents()  # M65: -22375 -> (cleared)
$ws.Cells.Item(65, 14).ClearContents()  # N65: -41240 -> (cleared)
$ws.Cells.Item(132, 8).Value = 23080010  # H132: 25003292 -> 23080010
$ws.Cells.Item(132, 9).Value = 30613458  # I132: 31916202 -> 30613458
$ws.Cells.Item(132, 10).Value = 8832.3125  # J132: 10462.154 -> 8832.3125
$ws.Cells.Item(132, 11).Value = 91840374  # K132: 95748606 -> 91840374
$ws.Cells.Item(132, 12).Value = 26496.9375  # L132: 31386.462 -> 26496.9375
$ws.Cells.Item(132, 13).Value = -91837844  # M132: -95746076 -> -91837844
$ws.Cells.Item(132, 14).Value = -31556.9375  # N132: -36446.462 -> -31556.9375
$ws.Cells.Item(140, 8).Value = 16714.5  # H140: 20000 -> 16714.5
$ws.Cells.Item(140, 10).Value = 16714.5  # J140: 20000 -> 16714.5
$ws.Cells.Item(140, 12).Value = 16714.5  # L140: 20000 -> 16714.5
$ws.Cells.Item(140, 14).Value = -27074.5  # N140: -30360 -> -27074.5
